# Applies odds updates to Sheet1 as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("S2").Value = 1.05
$ws.Range("F3").Value = 1.73
$ws.Range("G3").Value = 1.77
$ws.Range("N3").Value = 4.6
$ws.Range("P3").Value = 2.24
$ws.Range("Q3").Value = 1.74
$ws.Range("S3").Value = 2.84
$ws.Range("T3").Value = 1.74
$ws.Range("U3").Value = 2.22
$ws.Range("V3").Value = 1.22
$ws.Range("AA3").Value = 130
$ws.Range("AN3").Value = 9
$ws.Range("AO3").Value = 60
$ws.Range("F5").Value = 6.6
$ws.Range("G5").Value = 8
$ws.Range("I5").Value = 1.59
$ws.Range("K5").Value = 4.6
$ws.Range("N5").Value = 3.55
$ws.Range("Q5").Value = 1.98
$ws.Range("R5").Value = 1.33
$ws.Range("V5").Value = 2.68
$ws.Range("W5").Value = 1.14
$ws.Range("AE5").Value = 21
$ws.Range("AJ5").Value = 280
$ws.Range("AK5").Value = 140
$ws.Range("AN5").Value = 210
$ws.Range("AC6").Value = 10
$ws.Range("P7").Value = 2.82
$ws.Range("Q7").Value = 1.5
$ws.Range("R7").Value = 1.73
$ws.Range("S7").Value = 2
$ws.Range("T7").Value = 1.84
$ws.Range("F8").Value = 1.98
$ws.Range("G8").Value = 2.16
$ws.Range("H8").Value = 3.7
$ws.Range("J8").Value = 3.6
$ws.Range("K8").Value = 4.1
$ws.Range("Q8").Value = 1.7
$ws.Range("U8").Value = 2.18
$ws.Range("V8").Value = 1.3
$ws.Range("W8").Value = 1.86
$ws.Range("H9").Value = 2.28
$ws.Range("J9").Value = 3.5
$ws.Range("P9").Value = 2.06
$ws.Range("Q9").Value = 1.87
$ws.Range("W9").Value = 1.38
$ws.Range("AM9").Value = 80
$ws.Range("F10").Value = 3.8
$ws.Range("I10").Value = 2.04
$ws.Range("J10").Value = 3.35
$ws.Range("O10").Value = 1.23
$ws.Range("S10").Value = 2.6
$ws.Range("V10").Value = 1.96
$ws.Range("W10").Value = 1.26
$ws.Range("F12").Value = 2.36
$ws.Range("V12").Value = 1.44
$ws.Range("G13").Value = 1.69
$ws.Range("K13").Value = 6.2
$ws.Range("L13").Value = 1.41
$ws.Range("N13").Value = 2.34
$ws.Range("P13").Value = 1.58
$ws.Range("Q13").Value = 2.08
$ws.Range("W13").Value = 2.46
$ws.Range("J14").Value = 3.75
$ws.Range("P14").Value = 2.44
$ws.Range("Q14").Value = 1.56
$ws.Range("AC14").Value = 10.5
$ws.Range("AD14").Value = 14.5
$ws.Range("AE14").Value = 980
$ws.Range("AM14").Value = 65
$ws.Range("K15").Value = 5.3
$ws.Range("AL15").Value = 210
$ws.Range("G16").Value = 4.8
$ws.Range("K16").Value = 5
$ws.Range("R16").Value = 1.71
$ws.Range("S16").Value = 2.16
$ws.Range("W16").Value = 1.26
$ws.Range("X16").Value = 34
$ws.Range("AB16").Value = 30
$ws.Range("AC16").Value = 13.5
$ws.Range("AD16").Value = 13
$ws.Range("AE16").Value = 19.5
$ws.Range("AF16").Value = 48
$ws.Range("AG16").Value = 22
$ws.Range("AH16").Value = 19.5
$ws.Range("AI16").Value = 29
$ws.Range("AK16").Value = 48
$ws.Range("AL16").Value = 50
$ws.Range("AN16").Value = 36
$ws.Range("F17").Value = 2.56
$ws.Range("H17").Value = 2.48
$ws.Range("J17").Value = 3.85
$ws.Range("P17").Value = 2.42
$ws.Range("T17").Value = 1.54
$ws.Range("U17").Value = 2.52
$ws.Range("V17").Value = 1.58
$ws.Range("G18").Value = 1.82
$ws.Range("W18").Value = 2.2
$ws.Range("F19").Value = 2.6
$ws.Range("G19").Value = 2.62
$ws.Range("H19").Value = 3.05
$ws.Range("V19").Value = 1.47
$ws.Range("Z19").Value = 19
$ws.Range("AA19").Value = 50
$ws.Range("AC19").Value = 7.2
$ws.Range("AK19").Value = 28
$ws.Range("J20").Value = 3.6
$ws.Range("O20").Value = 1.47
$ws.Range("R20").Value = 1.25
$ws.Range("V20").Value = 1.22
$ws.Range("AB20").Value = 6.8
$ws.Range("F21").Value = 1.34
$ws.Range("H21").Value = 10.5
$ws.Range("I21").Value = 12
$ws.Range("S21").Value = 2.96
$ws.Range("T21").Value = 2.2
$ws.Range("V21").Value = 1.09
$ws.Range("Y21").Value = 34
$ws.Range("Z21").Value = 110
$ws.Range("AD21").Value = 40
$ws.Range("AE21").Value = 210
$ws.Range("AN21").Value = 5.9
$ws.Range("AO21").Value = 280
